# Amazon_data.xlsx modification:
# Insert two new leading columns (UserName / Password) in front of the
# existing Item_name / 65-inch TV data, matching shared-string ordering
# produced by the original author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing column A (and any content) two columns to the right,
# making room for the new "UserName"/"Password" columns at A:B. The old
# A1/A2 data (Item_name/65-inch TV) ends up in column C.
$ws.Columns("A:B").Insert()

# New header row
$ws.Range("A1").Value2 = "UserName"
$ws.Range("B1").Value2 = "Password"

# New data row
$ws.Range("A2").Value2 = "abc"
$ws.Range("B2").Value2 = "xyz"

# Match the highlighted-header fill that the original Item_name header
# (now in C1) already carries.
$ws.Range("A1:B1").Interior.Color = $ws.Range("C1").Interior.Color

# Leave the selection where the author's session left it.
[void]$ws.Range("E10").Select()
